$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-48
$data = @(
    @(6,7),
    @(8,8),
    @(9,9),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(11,11),
    @(9,9),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(5,6),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(8,8),
    @(8,8),
    @(6,6),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(8,8),
    @(7,7)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
